# Decisions.xlsx schema change (Fixed #198): add a DecisionKindId
# column, retire ColorSettingsId, and rework the seeded decision
# rows into the new "Оказать / Оказать срочно / Оказать планово"
# kinds with generic Begin/EndDateTime sentinel values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Header row (row 1): rename columns E / F / G
# ---------------------------------------------------------------
$ws.Range("E1").Value = "DecisionKindId"
$ws.Range("F1").Value = "BeginDateTime"
$ws.Range("G1").Value = "EndDateTime"

# ---------------------------------------------------------------
# BeginDateTime / EndDateTime (cols F / G): every row now shares
# the same generic "always active" window - serial 2 (1900-01-02)
# through serial 2958100 (9998-12-31), formatted as dates.
# ---------------------------------------------------------------
$ws.Range("F2").Value = 2
$ws.Range("F3").Value = 2
$ws.Range("F4").Value = 2

$ws.Range("G2").Value = 2958100
$ws.Range("G3").Value = 2958100
$ws.Range("G4").Value = 2958100

$ws.Range("F2").Copy()
$ws.Range("G2:G4").PasteSpecial(-4122)

# ---------------------------------------------------------------
# Name / ShortName text updates
# ---------------------------------------------------------------
$ws.Range("C2").Value = "Оказать"
$ws.Range("D2").Value = "Оказать"

$ws.Range("C3").Value = "Оказать срочно"
$ws.Range("D3").Value = "Оказать ср."

$ws.Range("C4").Value = "Оказать планово"
$ws.Range("D4").Value = "Оказать план."

# ---------------------------------------------------------------
# ParentId (col B) and the new DecisionKindId (col E) columns.
# Row 2 stays NULL/NULL; rows 3 and 4 now point at decision kind 1.
# Values are entered with a leading apostrophe so they are kept as
# text (matching the existing Id-style text columns), then the
# number format is normalized back to the sheet's default style.
# ---------------------------------------------------------------
$ws.Range("E2").Value = "'NULL"
$ws.Range("B3").Value = "'1"
$ws.Range("E3").Value = "'1"
$ws.Range("B4").Value = "'1"
$ws.Range("E4").Value = "'1"

$ws.Range("A2").Copy()
$ws.Range("E2").PasteSpecial(-4122)
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("E3").PasteSpecial(-4122)
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("E4").PasteSpecial(-4122)
